$d = $word.ActiveDocument

$titleText = "Play Deepsea Riches Free - Slot Game Review 2021"
$metaOldText = "Explore sunken pirate treasures in Deepsea Riches - read our review and play the game for free on your PC or mobile device."

# ---------------------------------------------------------------------------
# Locate the relevant paragraphs by content instead of trusting fixed
# indices: the document heading ("Heading 1" style) vs. the duplicated
# plain-text title near the bottom, plus the italic meta-description blurb
# that immediately follows it.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count

$headingIdx = -1
$dupTitleIdx = -1
$metaBlurbIdx = -1

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text.TrimEnd()
    if ($t -eq $titleText -and $p.Style.NameLocal -eq "Heading 1" -and $headingIdx -eq -1) {
        $headingIdx = $i
    }
    if ($t -eq $titleText -and $p.Style.NameLocal -ne "Heading 1") {
        $dupTitleIdx = $i
    }
    if ($t -eq $metaOldText) {
        $metaBlurbIdx = $i
    }
}

# ---------------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the title heading:
#    bold "Meta description" run followed by a plain run with the rest of
#    the sentence.
# ---------------------------------------------------------------------------
$headingPara = $d.Paragraphs($headingIdx)
$headingPara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs($headingIdx + 1)
$metaPara.Style = "Normal"

$metaLabel = "Meta description"
$metaRest  = ": Explore sunken pirate treasures in Deepsea Riches - read our review and play the game for free on your PC or mobile device."

$metaRange = $metaPara.Range
$metaRange.Text = $metaLabel + $metaRest

$metaStart = $metaPara.Range.Start
$boldRange = $d.Range($metaStart, $metaStart + $metaLabel.Length)
$boldRange.Font.Bold = 1

# ---------------------------------------------------------------------------
# 2. Near the end of the document, remove the duplicated bold title
#    paragraph entirely, and rewrite the following italic paragraph's text
#    into the DALLE image prompt (keeping its italic formatting).
#    (Indices above were computed before the insertion above; that
#    insertion happened earlier in the body, so every paragraph from
#    $headingIdx + 1 onward shifted down by one slot.)
# ---------------------------------------------------------------------------
$dupTitleIdx = $dupTitleIdx + 1
$metaBlurbIdx = $metaBlurbIdx + 1

$dupTitlePara = $d.Paragraphs($dupTitleIdx)
$dupTitlePara.Range.Delete()

# Deleting that whole paragraph (and its mark) shifts everything after it
# down by one slot as well.
$metaBlurbIdx = $metaBlurbIdx - 1
$metaBlurbPara = $d.Paragraphs($metaBlurbIdx)

$blurbStart = $metaBlurbPara.Range.Start
$blurbEnd   = $metaBlurbPara.Range.End
$bodyRange = $d.Range($blurbStart, $blurbEnd - 1)

$dallePrompt = 'Prompt for DALLE: Create a feature image for Deepsea Riches that is playful and adventurous. The image should be in cartoon style and feature a happy Maya warrior wearing glasses. The warrior should be diving underwater, surrounded by various symbols of the game, such as pirate artifacts, a treasure chest, and sea creatures like crabs and sharks. The colors should be bright and vibrant, with an emphasis on blues and greens to depict the underwater environment. The image should also include the title "Deepsea Riches" in bold letters at the top.'
$bodyRange.Text = $dallePrompt

Write-Output "Done. Paragraphs now: $($d.Paragraphs.Count)"
